$d = $word.ActiveDocument

# 1. "Challenges and Drawbacks" + " –" -> merge to "Challenges and Drawbacks –"
$d.Content.Find.Execute("Challenges and Drawbacks ", $true, $false, $false, $false, $false, $true, 1, $false, "Challenges and Drawbacks ", 2)

# 2. "So" with proofErr -> "So," (text only change here; bookmark handled separately below)
$d.Content.Find.Execute("So it does not know", $true, $false, $false, $false, $false, $true, 1, $false, "So, it does not know", 2)

# 3. "It adapts" merge
$d.Content.Find.Execute("It adapts", $true, $false, $false, $false, $false, $true, 1, $false, "It adapts", 2)

# 4. Adagrad paragraph merges
$d.Content.Find.Execute("the accumulated sum keeps growing", $true, $false, $false, $false, $false, $true, 1, $false, "the accumulated sum keeps growing", 2)
